$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 380

$ws1.Range("F3").Value = 837
$ws1.Range("G3").Value = "已售罄"

$ws1.Range("F4").Value = 283

$ws1.Range("F5").Value = 1037

$ws1.Range("F6").Value = 2413

$ws1.Range("F7").Value = 204

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 380

$ws4.Range("F3").Value = 837
$ws4.Range("G3").Value = "已售罄"

$ws4.Range("F4").Value = 283

$ws4.Range("F7").Value = 1037

$ws4.Range("F8").Value = 2413

$ws4.Range("F10").Value = 204
